# Edit sheet Card24 by admin
# Change the "card" value in column A (rows 2-12) from "2" to "24",
# keeping the cells stored as text (matching the sheet's existing
# text-formatted data in that column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$rng = $ws.Range("A2:A12")
$rng.NumberFormat = "@"
$rng.Value = "24"
